$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.447.13'
$ws.Range("E2").Value = '  -0.41%  '
$ws.Range("D3").Value = '1.896.72'
$ws.Range("E3").Value = '  +1.13%  '
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").Value = "'237.88"
$ws.Range("E5").Value = '  +1.09%  '
$ws.Range("E6").Value = '  -0.25%  '
$ws.Range("D7").Value = "'0.4911"
$ws.Range("E7").Value = '  +0.95%  '
$ws.Range("E8").Value = '  +0.89%  '
$ws.Range("D9").Value = "'0.06689"
$ws.Range("E9").Value = '  +0.72%  '
$ws.Range("D10").Value = '1.884.66'
$ws.Range("E10").Value = '  +0.44%  '
$ws.Range("D11").Value = "'16.88"
$ws.Range("E11").Value = '  +2.11%  '
$ws.Range("D12").Value = "'0.07328"
$ws.Range("E12").Value = '  +1.39%  '
$ws.Range("D13").Value = "'5.179"
$ws.Range("E13").Value = '  +3.94%  '
$ws.Range("D14").Value = "'87.22"
$ws.Range("E14").Value = '  -1.73%  '
$ws.Range("D15").Value = "'0.6642"
$ws.Range("E15").Value = '  +2.26%  '
$ws.Range("D16").Value = '30.419.52'
$ws.Range("E16").Value = '  -0.31%  '
$ws.Range("D17").Value = "'13.45"
$ws.Range("E17").Value = '  +4.13%  '
$ws.Range("D18").Value = "'0.000007821"
$ws.Range("E18").Value = '  +0.04%  '
$ws.Range("E19").Value = '  -0.20%  '
$ws.Range("D20").Value = '2.132.17'
$ws.Range("E20").Value = '  +0.22%  '
$ws.Range("D21").Value = "'5.349"
$ws.Range("E21").Value = '  +13.65%  '
$ws.Range("E22").Value = '  -0.28%  '
$ws.Range("D23").Value = "'193.19"
$ws.Range("E23").Value = '  +1.29%  '
$ws.Range("D24").Value = "'6.104"
$ws.Range("E24").Value = '  +1.04%  '
$ws.Range("D25").Value = "'9.466"
$ws.Range("E25").Value = '  +2.31%  '
$ws.Range("D26").Value = "'162.35"
$ws.Range("E26").Value = '  +2.64%  '
$ws.Range("D27").Value = "'18.19"
$ws.Range("E27").Value = '  -0.17%  '
$ws.Range("D28").Value = "'1.931"
$ws.Range("E28").Value = '  +6.15%  '
$ws.Range("D29").Value = "'1.470"
$ws.Range("E29").Value = '  +4.60%  '
$ws.Range("D30").Value = "'4.318"
$ws.Range("E30").Value = '  +2.27%  '
$ws.Range("D31").Value = "'0.09154"
$ws.Range("E31").Value = '  +1.95%  '
$ws.Range("D32").Value = "'4.050"
$ws.Range("E32").Value = '  +3.72%  '
$ws.Range("D33").Value = "'0.05165"
$ws.Range("E33").Value = '  +0.84%  '
$ws.Range("D34").Value = "'0.7366"
$ws.Range("E34").Value = '  +2.37%  '
$ws.Range("E35").Value = '  +2.52%  '
$ws.Range("E36").Value = '  +0.64%  '
$ws.Range("D37").Value = "'0.01807"
$ws.Range("E37").Value = '  +0.33%  '
$ws.Range("D38").Value = "'2.673"
$ws.Range("E38").Value = '  +0.75%  '
$ws.Range("D39").Value = "'0.9231"
$ws.Range("E39").Value = '  +0.78%  '
$ws.Range("D40").Value = "'2.037"
$ws.Range("E40").Value = '  +0.21%  '
$ws.Range("D41").Value = "'0.4377"
$ws.Range("E41").Value = '  +0.63%  '
$ws.Range("D42").Value = "'5.907"
$ws.Range("E42").Value = '  +3.76%  '
$ws.Range("E43").Value = '  +1.71%  '
$ws.Range("D44").Value = "'0.9940"
$ws.Range("E44").Value = '  -0.20%  '
$ws.Range("D45").Value = "'68.39"
$ws.Range("E45").Value = '  +20.34%  '
$ws.Range("D46").Value = "'0.1363"
$ws.Range("E46").Value = '  +2.86%  '
$ws.Range("D47").Value = "'7.572"
$ws.Range("E47").Value = '  +3.81%  '
$ws.Range("D48").Value = "'8.990"
$ws.Range("E48").Value = '  +4.30%  '
$ws.Range("D49").Value = "'34.89"
$ws.Range("E49").Value = '  +5.42%  '
$ws.Range("D50").Value = "'0.05847"
$ws.Range("E50").Value = '  +0.42%  '
$ws.Range("D51").Value = "'0.3913"
$ws.Range("E51").Value = '  -2.71%  '
